$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 0.5526037812232971
$ws.Range("B1").Value = 1.083119034767151
$ws.Range("C1").Value = 5.09970760345459
$ws.Range("D1").Value = 3.767892837524414
$ws.Range("E1").Value = 1.131724238395691
